$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 01:42"

# Update country stats that changed (new data pull + re-sort by Casos totales desc)
# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 8214494
$ws.Cells.Item(4,3).Value = 64308
$ws.Cells.Item(4,4).Value = 5311547
$ws.Cells.Item(4,5).Value = 2680259
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 845
$ws.Cells.Item(4,8).Value = 222688

# Row 6: Brasil
$ws.Cells.Item(6,2).Value = 5170996
$ws.Cells.Item(6,3).Value = 29498
$ws.Cells.Item(6,4).Value = 4599446
$ws.Cells.Item(6,5).Value = 419037
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 734
$ws.Cells.Item(6,8).Value = 152513

# Row 10: Colombia
$ws.Cells.Item(10,2).Value = 936982
$ws.Cells.Item(10,3).Value = 6823
$ws.Cells.Item(10,4).Value = 826831
$ws.Cells.Item(10,5).Value = 81694
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 151
$ws.Cells.Item(10,8).Value = 28457

# Row 30: Canada
$ws.Cells.Item(30,2).Value = 191730
$ws.Cells.Item(30,3).Value = 2343
$ws.Cells.Item(30,4).Value = 161490
$ws.Cells.Item(30,5).Value = 20541
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 35
$ws.Cells.Item(30,8).Value = 9699

# Row 36: Chequia
$ws.Cells.Item(36,2).Value = 149010
$ws.Cells.Item(36,3).Value = 9720
$ws.Cells.Item(36,4).Value = 63350
$ws.Cells.Item(36,5).Value = 84430
$ws.Cells.Item(36,6).Value = 0
$ws.Cells.Item(36,7).Value = 58
$ws.Cells.Item(36,8).Value = 1230

# Row 41: Republica Dominicana
$ws.Cells.Item(41,2).Value = 120066
$ws.Cells.Item(41,3).Value = 404
$ws.Cells.Item(41,4).Value = 96152
$ws.Cells.Item(41,5).Value = 21725
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 3
$ws.Cells.Item(41,8).Value = 2189

# Row 61: Nigeria
$ws.Cells.Item(61,2).Value = 60982
$ws.Cells.Item(61,3).Value = 148
$ws.Cells.Item(61,4).Value = 52194
$ws.Cells.Item(61,5).Value = 7672
$ws.Cells.Item(61,6).Value = 0
$ws.Cells.Item(61,7).Value = 0
$ws.Cells.Item(61,8).Value = 1116

# Row 67: Paraguay
$ws.Cells.Item(67,2).Value = 52596
$ws.Cells.Item(67,3).Value = 751
$ws.Cells.Item(67,4).Value = 34427
$ws.Cells.Item(67,5).Value = 17019
$ws.Cells.Item(67,6).Value = 0
$ws.Cells.Item(67,7).Value = 19
$ws.Cells.Item(67,8).Value = 1150

# Row 85: Australia
$ws.Cells.Item(85,2).Value = 27362
$ws.Cells.Item(85,3).Value = 21
$ws.Cells.Item(85,4).Value = 25061
$ws.Cells.Item(85,5).Value = 1397
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = 904

# Row 96: Noruega
$ws.Cells.Item(96,2).Value = 16137
$ws.Cells.Item(96,3).Value = 184
$ws.Cells.Item(96,4).Value = 11863
$ws.Cells.Item(96,5).Value = 3996
$ws.Cells.Item(96,6).Value = 0
$ws.Cells.Item(96,7).Value = 1
$ws.Cells.Item(96,8).Value = 278

# Row 97: Zambia
$ws.Cells.Item(97,2).Value = 15659
$ws.Cells.Item(97,3).Value = 43
$ws.Cells.Item(97,4).Value = 14899
$ws.Cells.Item(97,5).Value = 414
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 1
$ws.Cells.Item(97,8).Value = 346

# Row 113: Haiti
$ws.Cells.Item(113,2).Value = 8925
$ws.Cells.Item(113,3).Value = 17
$ws.Cells.Item(113,4).Value = 7182
$ws.Cells.Item(113,5).Value = 1512
$ws.Cells.Item(113,6).Value = 0
$ws.Cells.Item(113,7).Value = 0
$ws.Cells.Item(113,8).Value = 231

# Row 119: Guadalupe -> Angola (re-sorted)
$ws.Cells.Item(119,1).Value = "Angola"
$ws.Cells.Item(119,2).Value = 7096
$ws.Cells.Item(119,3).Value = 250
$ws.Cells.Item(119,4).Value = 2928
$ws.Cells.Item(119,5).Value = 3940
$ws.Cells.Item(119,6).Value = 0
$ws.Cells.Item(119,7).Value = 1
$ws.Cells.Item(119,8).Value = 228

# Row 120: Angola -> Guadalupe (re-sorted)
$ws.Cells.Item(120,1).Value = "Guadalupe"
$ws.Cells.Item(120,2).Value = 6908
$ws.Cells.Item(120,3).Value = 0
$ws.Cells.Item(120,4).Value = 2199
$ws.Cells.Item(120,5).Value = 4613
$ws.Cells.Item(120,6).Value = 0
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 96

# Row 132: Surinam
$ws.Cells.Item(132,2).Value = 5094
$ws.Cells.Item(132,3).Value = 11
$ws.Cells.Item(132,4).Value = 4907
$ws.Cells.Item(132,5).Value = 78
$ws.Cells.Item(132,6).Value = 0
$ws.Cells.Item(132,7).Value = 1
$ws.Cells.Item(132,8).Value = 109

# Row 155: Guinea-Bisau -> Uruguay (re-sorted)
$ws.Cells.Item(155,1).Value = "Uruguay"
$ws.Cells.Item(155,2).Value = 2417
$ws.Cells.Item(155,3).Value = 29
$ws.Cells.Item(155,4).Value = 2025
$ws.Cells.Item(155,5).Value = 341
$ws.Cells.Item(155,6).Value = 0
$ws.Cells.Item(155,7).Value = 0
$ws.Cells.Item(155,8).Value = 51

# Row 156: Uruguay -> Guinea-Bisau (re-sorted)
$ws.Cells.Item(156,1).Value = "Guinea-Bisau"
$ws.Cells.Item(156,2).Value = 2389
$ws.Cells.Item(156,3).Value = 0
$ws.Cells.Item(156,4).Value = 1782
$ws.Cells.Item(156,5).Value = 566
$ws.Cells.Item(156,6).Value = 0
$ws.Cells.Item(156,7).Value = 0
$ws.Cells.Item(156,8).Value = 41

# Row 157: Sierra Leona -> Burkina Faso (re-sorted)
$ws.Cells.Item(157,1).Value = "Burkina Faso"
$ws.Cells.Item(157,2).Value = 2335
$ws.Cells.Item(157,3).Value = 30
$ws.Cells.Item(157,4).Value = 1645
$ws.Cells.Item(157,5).Value = 625
$ws.Cells.Item(157,6).Value = 0
$ws.Cells.Item(157,7).Value = 2
$ws.Cells.Item(157,8).Value = 65

# Row 158: Burkina Faso -> Sierra Leona (re-sorted)
$ws.Cells.Item(158,1).Value = "Sierra Leona"
$ws.Cells.Item(158,2).Value = 2323
$ws.Cells.Item(158,3).Value = 8
$ws.Cells.Item(158,4).Value = 1746
$ws.Cells.Item(158,5).Value = 504
$ws.Cells.Item(158,6).Value = 0
$ws.Cells.Item(158,7).Value = 0
$ws.Cells.Item(158,8).Value = 73

# Row 171: San Martin (Parte Holandesa)
$ws.Cells.Item(171,2).Value = 737
$ws.Cells.Item(171,3).Value = 8
$ws.Cells.Item(171,4).Value = 657
$ws.Cells.Item(171,5).Value = 58
$ws.Cells.Item(171,6).Value = 0
$ws.Cells.Item(171,7).Value = 0
$ws.Cells.Item(171,8).Value = 22

# Row 209: Granada -> Santa Sede (re-sorted)
$ws.Cells.Item(209,1).Value = "Santa Sede"
$ws.Cells.Item(209,2).Value = 26
$ws.Cells.Item(209,3).Value = 7
$ws.Cells.Item(209,4).Value = 12
$ws.Cells.Item(209,5).Value = 14
$ws.Cells.Item(209,6).Value = 0
$ws.Cells.Item(209,7).Value = 0
$ws.Cells.Item(209,8).Value = 0

# Row 210: Laos -> Granada (re-sorted)
$ws.Cells.Item(210,1).Value = "Granada"
$ws.Cells.Item(210,2).Value = 25
$ws.Cells.Item(210,3).Value = 0
$ws.Cells.Item(210,4).Value = 24
$ws.Cells.Item(210,5).Value = 1
$ws.Cells.Item(210,6).Value = 0
$ws.Cells.Item(210,7).Value = 0
$ws.Cells.Item(210,8).Value = 0

# Row 211: Santa Sede -> Laos (re-sorted)
$ws.Cells.Item(211,1).Value = "Laos"
$ws.Cells.Item(211,2).Value = 23
$ws.Cells.Item(211,3).Value = 0
$ws.Cells.Item(211,4).Value = 22
$ws.Cells.Item(211,5).Value = 1
$ws.Cells.Item(211,6).Value = 0
$ws.Cells.Item(211,7).Value = 0
$ws.Cells.Item(211,8).Value = 0

